$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextCell 2 4 "71.953.70"
$ws.Cells.Item(2, 5).Value = "  +4.98%  "
Set-TextCell 3 4 "2.632.91"
$ws.Cells.Item(3, 5).Value = "  +4.80%  "
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
Set-TextCell 5 4 "608.66"
$ws.Cells.Item(5, 5).Value = "  +2.91%  "
Set-TextCell 6 4 "179.70"
$ws.Cells.Item(6, 5).Value = "  +2.95%  "
$ws.Cells.Item(7, 5).Value = "  -0.04%  "
Set-TextCell 8 4 "0.525"
$ws.Cells.Item(8, 5).Value = "  +1.82%  "
Set-TextCell 9 4 "2.630.78"
$ws.Cells.Item(9, 5).Value = "  +4.68%  "
Set-TextCell 10 4 "0.169"
$ws.Cells.Item(10, 5).Value = "  +14.67%  "
$ws.Cells.Item(11, 5).Value = "  +1.05%  "
Set-TextCell 12 4 "0.349"
$ws.Cells.Item(12, 5).Value = "  +3.99%  "
$ws.Cells.Item(13, 5).Value = "  +0.39%  "
$ws.Cells.Item(14, 5).Value = "  +10.18%  "
Set-TextCell 15 4 "3.084.19"
$ws.Cells.Item(15, 5).Value = "  +3.56%  "
Set-TextCell 16 4 "71.840.42"
$ws.Cells.Item(16, 5).Value = "  +4.92%  "
Set-TextCell 17 4 "26.54"
$ws.Cells.Item(17, 5).Value = "  +2.97%  "
Set-TextCell 18 4 "2.639.79"
$ws.Cells.Item(18, 5).Value = "  +4.90%  "
Set-TextCell 19 4 "8.05"
$ws.Cells.Item(19, 5).Value = "  +7.06%  "
Set-TextCell 20 4 "383.52"
$ws.Cells.Item(20, 5).Value = "  +6.32%  "
Set-TextCell 21 4 "11.51"
$ws.Cells.Item(21, 5).Value = "  +5.74%  "
$ws.Cells.Item(22, 5).Value = "  +3.23%  "
$ws.Cells.Item(23, 5).Value = "  +21.67%  "
Set-TextCell 24 4 "72.76"
$ws.Cells.Item(24, 5).Value = "  +3.84%  "
$ws.Cells.Item(25, 5).Value = "  +6.99%  "
$ws.Cells.Item(26, 5).Value = "  +0.04%  "
Set-TextCell 27 4 "9.96"
$ws.Cells.Item(27, 5).Value = "  +12.51%  "
Set-TextCell 28 4 "2.767.13"
$ws.Cells.Item(28, 5).Value = "  +4.67%  "
$ws.Cells.Item(29, 5).Value = "  +1.27%  "
Set-TextCell 30 4 "0.0₃0968"
$ws.Cells.Item(30, 5).Value = "  +11.04%  "
Set-TextCell 31 4 "546.19"
$ws.Cells.Item(31, 5).Value = "  +7.45%  "
Set-TextCell 32 4 "8.07"
$ws.Cells.Item(32, 5).Value = "  +4.58%  "
$ws.Cells.Item(33, 5).Value = "  +9.29%  "
$ws.Cells.Item(34, 5).Value = "  +3.83%  "
Set-TextCell 35 4 "1.00"
$ws.Cells.Item(35, 5).Value = "  -0.06%  "
Set-TextCell 36 4 "166.11"
$ws.Cells.Item(36, 5).Value = "  +2.69%  "
Set-TextCell 37 4 "19.26"
$ws.Cells.Item(37, 5).Value = "  +3.84%  "
Set-TextCell 38 4 "0.115"
$ws.Cells.Item(38, 5).Value = "  -2.22%  "
Set-TextCell 39 4 "19.11"
$ws.Cells.Item(39, 5).Value = "  +2.58%  "
Set-TextCell 41 4 "1.86"
$ws.Cells.Item(41, 5).Value = "  +9.37%  "
$ws.Cells.Item(42, 5).Value = "  +0.13%  "
$ws.Cells.Item(43, 5).Value = "  +11.86%  "
Set-TextCell 44 4 "5.04"
$ws.Cells.Item(44, 5).Value = "  +6.37%  "
Set-TextCell 45 4 "0.333"
$ws.Cells.Item(45, 5).Value = "  +4.14%  "
Set-TextCell 46 4 "39.62"
$ws.Cells.Item(46, 5).Value = "  +1.80%  "
Set-TextCell 47 4 "151.02"
$ws.Cells.Item(47, 5).Value = "  +0.61%  "
Set-TextCell 48 4 "3.66"
$ws.Cells.Item(48, 5).Value = "  +3.12%  "
$ws.Cells.Item(49, 5).Value = "  +5.26%  "
Set-TextCell 50 4 "1.70"
Set-TextCell 51 4 "0.0₆0265"
$ws.Cells.Item(51, 5).Value = "  +6.45%  "
